$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 ("Update Manual" task on 2013-06-07): the 4h effort is split into
# 2h of regular effort (B8) plus 2h of additional effort (C8).
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 2

# Row 9 (new): another session on the same day continuing
# "Implementation of mutexes", 2.25h of effort.
# Copy A8's formatting (date number format) down to A9 first, then set
# the actual date value, so the new cell reuses the existing date style
# instead of creating a duplicate.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial()
$ws.Range("A9").Value = $ws.Range("A8").Value2
$ws.Range("B9").Value = 2.25
$ws.Range("D9").Value = "Implementation of mutexes"

# Row 10 (new): entry for 2013-06-10, mutexes basically done.
$ws.Range("A8").Copy()
$ws.Range("A10").PasteSpecial()
$ws.Range("A10").Value = 41435
$ws.Range("B10").Value = 2
$ws.Range("D10").Value = "Implementation of mutexes: Basically done. No test case implemented yet, no testing done yet"

# Move the selection to reflect where the user ended up after editing.
$ws.Range("E10").Select()
